# Add season-record columns (Wins / Losses / Ties) to the team sheet.
#
# The sheet currently ends at column AC (Unnamed: 28). We extend it with
# three new columns: AD = Wins, AE = Losses, AF = Ties, matching the
# formatting of the existing header row and filling every player row
# (2-50) with the team's 1996 season record (92 wins, 70 losses, 0 ties).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reuse the existing header formatting (bold, bordered, centered/top
# aligned) from the last current header cell (AC1) for the new headers.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$wins = 92
$losses = 70
$ties = 0

$lastRow = 50
for ($row = 2; $row -le $lastRow; $row++) {
    $ws.Cells.Item($row, 30).Value = $wins
    $ws.Cells.Item($row, 31).Value = $losses
    $ws.Cells.Item($row, 32).Value = $ties
}
